# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet at the very front of the workbook.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, and shrink the stored value
#    down from the full scorecard URL to just the bare match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet before everything else ---------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Style = "Normal"
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4712"
$playerInfo.Range("B2").Value = "Craig Overton"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# Match the page-margin convention used by the rest of the workbook's
# sheets (0.75in sides, 1in top/bottom, 0.5in header/footer).
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

# Re-fetch these by name now that the sheet collection has shifted
# (worksheet variables captured before the Add() above track *position*,
# not the sheet object, so they'd otherwise now point at "Player Info").
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2. ODI Batting: MATCH_CARD_LINK (col D) -> MATCH_CODE --------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

# Row 2..8, in order -> bare match code (was the full scorecard URL).
$matchCodes = @("4169", "4472", "4473", "4476", "4609", "4613", "4618")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $cell = $battingSheet.Cells.Item($i + 2, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# --- 3. ODI Bowling: MATCH_CARD_LINK (col B) -> MATCH_CODE --------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $cell = $bowlingSheet.Cells.Item($i + 2, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}
